# Append the 09/19/2025 profit-run row to the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds the date as literal text (like the other rows in this
# sheet), not a real date serial. Force text typing via NumberFormat so
# Excel doesn't auto-convert the "mm/dd/yyyy"-looking string into a date,
# then clear the format again so the cell keeps the sheet's default style
# (matching how the existing text-date cells are styled).
$ws.Range("A33").NumberFormat = "@"
$ws.Range("A33").Value = "09/19/2025"
$ws.Range("A33").ClearFormats()

$ws.Range("B33").Value = 15979.25
